# fixed typos on Slide 19 (and related Prerequisites slide cleanup)
$p = $ppt.ActivePresentation

# --- Slide 19 ("Single Stage Iterative Workflow"): the callout box
# mistakenly referenced Stage_3 / Stage_4 (copy/pasted from the "Multi
# Stage" slide) instead of Stage_1 / Stage_1, which is what this
# single-stage slide actually uses everywhere else.
$s19 = $p.Slides.Item(19)
$callout = $s19.Shapes.Item("Rounded Rectangular Callout 8")
$tr19 = $callout.TextFrame.TextRange

# "Stage_3.Input_1 and Stage_4.Output_1 should have IDENTICAL number of files"
#  12345678 9      17  21      29...
$tr19.Characters(1, 8).Text = "Stage_1."
$tr19.Characters(9, 8).Text = "Input_1 "
$tr19.Characters(17, 4).Text = "and "
$tr19.Characters(21, 8).Text = "Stage_1."

# --- Slide 2 ("Prerequisites"): collapse the word-by-word runs back into
# single runs (no per-word formatting differences remain).
$s2 = $p.Slides.Item(2)
$content = $s2.Shapes.Item("Content Placeholder 2")
$tr2 = $content.TextFrame.TextRange

$f1 = $tr2.Find("Python2 is not compatible with the external mapper")
$tr2.Characters($f1.Start, $f1.Length).Text = "Python2 is not compatible with the external mapper"

$f2 = $tr2.Find("Requires Pegasus to run output Pegasus DAX")
$tr2.Characters($f2.Start, $f2.Length).Text = "Requires Pegasus to run output Pegasus DAX"
